# Auto-generated edit script applying diff changes to cryptos worksheet
# Updates coin price (D) and volume/1h change (E) columns, plus
# re-orders a few coin rows (B/C/D/E) to match the refreshed ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '44.266.60'

# Row 3
$ws.Cells.Item(3, 4).Value = '2.261.36'
$ws.Cells.Item(3, 5).Value = '  -0.33%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.16%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '306.34'
$ws.Cells.Item(5, 5).Value = '  -5.98%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '99.08'
$ws.Cells.Item(6, 5).Value = '  -4.77%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.576'
$ws.Cells.Item(7, 5).Value = '  -2.28%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.10%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.542'
$ws.Cells.Item(9, 5).Value = '  -4.78%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '36.15'
$ws.Cells.Item(10, 5).Value = '  -6.73%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0825'
$ws.Cells.Item(11, 5).Value = '  -2.56%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '7.41'
$ws.Cells.Item(12, 5).Value = '  -6.22%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.104'
$ws.Cells.Item(13, 5).Value = '  -2.67%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.605.13'
$ws.Cells.Item(14, 5).Value = '  -0.35%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.848'
$ws.Cells.Item(15, 5).Value = '  -3.68%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.256.10'
$ws.Cells.Item(16, 5).Value = '  -0.55%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '14.01'
$ws.Cells.Item(17, 5).Value = '  -3.67%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '44.129.52'
$ws.Cells.Item(18, 5).Value = '  -0.11%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '13.10'
$ws.Cells.Item(19, 5).Value = '  -8.29%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0981'
$ws.Cells.Item(20, 5).Value = '  -2.47%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.39'
$ws.Cells.Item(21, 5).Value = '  -2.65%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '65.68'
$ws.Cells.Item(22, 5).Value = '  -1.21%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '242.07'
$ws.Cells.Item(23, 5).Value = '  +0.98%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.94'
$ws.Cells.Item(24, 5).Value = '  -9.79%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.01'
$ws.Cells.Item(25, 5).Value = '  -8.48%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.38%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '10.24'
$ws.Cells.Item(27, 5).Value = '  -0.61%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '37.86'
$ws.Cells.Item(28, 5).Value = '  -1.87%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.12'
$ws.Cells.Item(29, 5).Value = '  -3.76%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '6.11'
$ws.Cells.Item(30, 5).Value = '  -5.75%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '20.29'
$ws.Cells.Item(31, 5).Value = '  -1.39%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Monero'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '157.02'
$ws.Cells.Item(32, 5).Value = '  -3.50%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.0840'
$ws.Cells.Item(33, 5).Value = '  -4.65%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '3.47'
$ws.Cells.Item(34, 5).Value = '  +11.03%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.67'
$ws.Cells.Item(35, 5).Value = '  -1.96%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.93'
$ws.Cells.Item(36, 5).Value = '  -2.87%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Stellar'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.119'
$ws.Cells.Item(37, 5).Value = '  -1.96%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Kaspa'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.108'
$ws.Cells.Item(38, 5).Value = '  -8.37%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '16.04'
$ws.Cells.Item(39, 5).Value = '  +1.42%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.51'
$ws.Cells.Item(40, 5).Value = '  -10.50%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.95'
$ws.Cells.Item(41, 5).Value = '  -9.91%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.0310'
$ws.Cells.Item(42, 5).Value = '  -5.13%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.03%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'BitcoinSV'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '89.70'
$ws.Cells.Item(44, 5).Value = '  +5.79%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Maker'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(45, 4).Value = '1.731.85'
$ws.Cells.Item(45, 5).Value = '  -3.20%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.195'
$ws.Cells.Item(46, 5).Value = '  -5.32%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '5.19'
$ws.Cells.Item(47, 5).Value = '  -3.93%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '102.57'
$ws.Cells.Item(48, 5).Value = '  -3.19%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '71.74'
$ws.Cells.Item(49, 5).Value = '  -6.91%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '56.33'
$ws.Cells.Item(50, 5).Value = '  -6.72%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '8.29'
$ws.Cells.Item(51, 5).Value = '  -3.83%  '

Write-Output "cryptos list updated"